$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header columns E1/F1 (shared strings will be re-indexed automatically,
# dropping the now-unused old strings and appending the new ones)
$ws.Range("E1").Value = "percent_impared_water_island"
$ws.Range("F1").Value = "percent_impared_state"

# Resize / add columns C, E, F
$ws.Columns.Item(3).ColumnWidth = 28.0
$ws.Columns.Item(5).ColumnWidth = 29.833333333333332
$ws.Columns.Item(6).ColumnWidth = 20.333333333333332

# Update the active selection
$ws.Activate()
$ws.Range("H9").Select()
